# Updated cryptos list values (price + 1h volume change) per the target diff.
# Some rows also swap the coin Name (B) / Link (C) columns to reflect re-ranked order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (e.g. '0.630', '1.00', '65.823.57')
# that Excel would otherwise auto-convert to a Number and mangle (dropping
# trailing zeros / failing to parse multi-dot values). Force Text format on
# each target cell right before the write, then restore the cell's original
# (default/"Normal") style so no visible formatting changes are introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

# Row 2
Set-TextValue $ws.Range('D2') '65.642.92'
$ws.Range('E2').Value = '  -2.66%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.268.10'
$ws.Range('E3').Value = '  -1.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '571.93'
$ws.Range('E5').Value = '  -0.78%  '

# Row 6
Set-TextValue $ws.Range('D6') '176.55'
$ws.Range('E6').Value = '  -4.96%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.630'
$ws.Range('E7').Value = '  +4.45%  '

# Row 8
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.125'
$ws.Range('E9').Value = '  -3.22%  '

# Row 10
Set-TextValue $ws.Range('D10') '6.70'
$ws.Range('E10').Value = '  +0.44%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.398'
$ws.Range('E11').Value = '  -3.06%  '

# Row 12
Set-TextValue $ws.Range('D12') '3.840.52'
$ws.Range('E12').Value = '  -1.20%  '

# Row 13
$ws.Range('E13').Value = '  -3.90%  '

# Row 14
$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D14') '65.761.19'
$ws.Range('E14').Value = '  -2.77%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D15') '26.40'
$ws.Range('E15').Value = '  -3.82%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D16') '3.291.66'
$ws.Range('E16').Value = '  -1.20%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.0000162'
$ws.Range('E17').Value = '  -3.06%  '

# Row 18
Set-TextValue $ws.Range('D18') '434.73'
$ws.Range('E18').Value = '  -2.04%  '

# Row 19
Set-TextValue $ws.Range('D19') '5.55'
$ws.Range('E19').Value = '  -2.33%  '

# Row 20
Set-TextValue $ws.Range('D20') '13.12'
$ws.Range('E20').Value = '  -3.33%  '

# Row 21
Set-TextValue $ws.Range('D21') '7.36'
$ws.Range('E21').Value = '  -4.99%  '

# Row 22
Set-TextValue $ws.Range('D22') '72.36'
$ws.Range('E22').Value = '  -2.30%  '

# Row 23
Set-TextValue $ws.Range('D23') '0.999'
$ws.Range('E23').Value = '  +0.07%  '

# Row 24
Set-TextValue $ws.Range('D24') '3.429.25'
$ws.Range('E24').Value = '  -0.77%  '

# Row 25
Set-TextValue $ws.Range('D25') '0.506'
$ws.Range('E25').Value = '  -1.90%  '

# Row 26
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D26') '0.194'
$ws.Range('E26').Value = '  +3.44%  '

# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D27') '0.0000112'
$ws.Range('E27').Value = '  -5.39%  '

# Row 28
Set-TextValue $ws.Range('D28') '8.82'
$ws.Range('E28').Value = '  -2.46%  '

# Row 29
Set-TextValue $ws.Range('D29') '0.999'
$ws.Range('E29').Value = '  -0.03%  '

# Row 31
Set-TextValue $ws.Range('D31') '22.18'
$ws.Range('E31').Value = '  -3.27%  '

# Row 33
Set-TextValue $ws.Range('D33') '5.12'
$ws.Range('E33').Value = '  -3.94%  '

# Row 34
Set-TextValue $ws.Range('D34') '6.58'
$ws.Range('E34').Value = '  -3.42%  '

# Row 35
$ws.Range('E35').Value = '  -5.38%  '

# Row 36
Set-TextValue $ws.Range('D36') '158.67'
$ws.Range('E36').Value = '  -2.55%  '

# Row 37
$ws.Range('E37').Value = '  -5.09%  '

# Row 38
Set-TextValue $ws.Range('D38') '26.60'
$ws.Range('E38').Value = '  -2.31%  '

# Row 39
Set-TextValue $ws.Range('D39') '1.77'
$ws.Range('E39').Value = '  -4.31%  '

# Row 40
Set-TextValue $ws.Range('D40') '2.771.12'
$ws.Range('E40').Value = '  +0.67%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.773'
$ws.Range('E41').Value = '  -2.19%  '

# Row 42
Set-TextValue $ws.Range('D42') '4.30'
$ws.Range('E42').Value = '  -3.82%  '

# Row 43
Set-TextValue $ws.Range('D43') '40.18'
$ws.Range('E43').Value = '  -0.01%  '

# Row 44
$ws.Range('E44').Value = '  -3.62%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.0654'
$ws.Range('E45').Value = '  -2.64%  '

# Row 46
$ws.Range('E46').Value = '  -5.34%  '

# Row 47
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D47') '316.69'
$ws.Range('E47').Value = '  -3.04%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D48') '23.25'
$ws.Range('E48').Value = '  -6.40%  '

# Row 49
Set-TextValue $ws.Range('D49') '0.0267'
$ws.Range('E49').Value = '  -2.41%  '

# Row 50
$ws.Range('E50').Value = '  +2.38%  '

# Row 51
Set-TextValue $ws.Range('D51') '1.00'
$ws.Range('E51').Value = '  +0.00%  '
